# Slide 24 ("Entity Question"), Content Placeholder 4: the code listing
# paragraph "cascade = {...} , orphanRemoval = true)" currently ends with a
# single run whose text is "true)". Split that run into two runs - "true"
# and ")" - matching the author's edit (keystroke landed between "true" and
# the closing paren), while preserving the existing Courier New formatting.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(24)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$fullText = $tr.Text
$target = "true)"
$startPos = $fullText.IndexOf($target) + 1

if ($startPos -le 0) {
    throw "Could not locate 'true)' run in the shape's text"
}

# Re-assigning the sub-range's own text (same length) forces the host to
# split the underlying run in two at the boundary we choose, while each
# half keeps inheriting the original run's formatting (Courier New, lang).
$trueRange = $tr.Characters($startPos, 4)
$trueRange.Text = "true"

$parenRange = $tr.Characters($startPos + 4, 1)
$parenRange.Text = ")"
